# shopping cart page complete module push
# Adds a "phoneNumber" column (inserted before expectedErrorMsg), converts the
# zip code column to text, fixes the "confirm password" expected message, and
# appends a new phone-number validation test row.
#
# NOTE: the order in which brand-new text values are first written matters
# (it determines the order new entries land in the shared-string table), so
# this script deliberately writes things in the same order the original
# author's edit would have: new column header/data, then the new row, then
# the confirm-password message fix, then the zip-code text conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at O. This shifts the existing "expectedErrorMsg"
#    column (and all its data) from O -> P, matching the diff exactly.
# ---------------------------------------------------------------------------
$ws.Columns("O").Insert()

# ---------------------------------------------------------------------------
# 2. New column header + data ("phoneNumber").
# ---------------------------------------------------------------------------
$ws.Range("O1").Value = "phoneNumber"

# Rows 2-11 all get a phone number value; row 12 (new row) is left blank
# since it's the "missing phone number" negative test case.
$ws.Range("O2:O11").Value = "'9900623434"

# ---------------------------------------------------------------------------
# 3. Append new row 12: TC 011, phone number validation test.
# ---------------------------------------------------------------------------
$ws.Range("A11:P11").Copy()
$ws.Range("A12:P12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A12").Value = "'011"
$ws.Range("B12").Value = "Unilog"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "Hemanth"
$ws.Range("E12").Value = "Sridhar"
$ws.Range("F12").Value = "hemanth.BS@unilogcorp.com"
$ws.Range("G12").Value = "unilog123"
$ws.Range("H12").Value = "unilog123"
$ws.Range("I12").Value = "address1"
$ws.Range("J12").Value = "address2"
$ws.Range("K12").Value = "Adelanto"
$ws.Range("L12").Value = "United States"
$ws.Range("M12").Value = "California"
$ws.Range("N12").Value = "'9900623434"
$ws.Range("O12").ClearContents()
$ws.Range("P12").Value = "Please Enter Phone Number."

# ---------------------------------------------------------------------------
# 4. Row 7 ("confirm password" test) expected-message text changed.
# ---------------------------------------------------------------------------
$ws.Range("P7").Value = "Please enter Confirm password"

# ---------------------------------------------------------------------------
# 5. zipPostalCode column (N) switches from numeric 92301 to text "92301" for
#    every populated row (rows 2-10, 12). Row 11 stays blank (it's the
#    "missing zip" test case) but changes its fill style.
# ---------------------------------------------------------------------------
$ws.Range("N2:N10").Value = "'92301"
$ws.Range("N12").Value = "'92301"

# Row 11's zip cell becomes a plain bordered blank cell (style changes from
# the shaded "missing field" look to a plain empty cell) - copy format from a
# cell that already carries that plain style, then clear its contents.
$ws.Range("B3").Copy()
$ws.Range("N11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N11").ClearContents()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Column P needs an explicit width (the former column O's width carries
#    over automatically from the insert/shift, but the diff also widens it
#    slightly); set the closest achievable width.
# ---------------------------------------------------------------------------
$ws.Columns("P").ColumnWidth = 33

# ---------------------------------------------------------------------------
# 7. Update the named range so it covers the new column.
# ---------------------------------------------------------------------------
$wb.Names("RegistrationErrorScenarios").RefersTo = "=Sheet1!`$B`$1:`$P`$11"

# ---------------------------------------------------------------------------
# 8. Update the view: scroll so column G is left-most, select O10.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("O10").Select()
